$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '27.356.77'
$ws.Range("E2").Value = '  +9.20%  '
Set-TextValue "D3" '1.788.02'
$ws.Range("E3").Value = '  +6.65%  '
Set-TextValue "D4" '1.002'
$ws.Range("E4").Value = '  +0.12%  '
Set-TextValue "D5" '338.58'
$ws.Range("E5").Value = '  +2.68%  '
Set-TextValue "D6" '0.9991'
$ws.Range("E6").Value = '  -0.01%  '
Set-TextValue "D7" '0.3795'
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue "D8" '0.3500'
$ws.Range("E8").Value = '  +7.68%  '
$ws.Range("B9").Value = 'OKB'
$ws.Range("C9").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue "D9" '49.67'
$ws.Range("E9").Value = '  +4.67%  '
Set-TextValue "D10" '1.224'
$ws.Range("E10").Value = '  +6.34%  '
Set-TextValue "D11" '0.07637'
$ws.Range("E11").Value = '  +5.14%  '
Set-TextValue "D12" '1.0000'
$ws.Range("E12").Value = '  +0.06%  '
Set-TextValue "D13" '6.633'
$ws.Range("E13").Value = '  +8.73%  '
Set-TextValue "D14" '21.58'
$ws.Range("E14").Value = '  +9.33%  '
Set-TextValue "D15" '7.216'
$ws.Range("E15").Value = '  +8.11%  '
Set-TextValue "D16" '1.788.09'
$ws.Range("E16").Value = '  +6.46%  '
Set-TextValue "D17" '0.00001117'
$ws.Range("E17").Value = '  +6.00%  '
Set-TextValue "D18" '0.06786'
$ws.Range("E18").Value = '  +3.66%  '
Set-TextValue "D19" '85.46'
$ws.Range("E19").Value = '  +8.05%  '
Set-TextValue "D20" '0.9993'
$ws.Range("E20").Value = '  +0.06%  '
Set-TextValue "D21" '17.66'
$ws.Range("E21").Value = '  +11.16%  '
Set-TextValue "D22" '6.424'
$ws.Range("E22").Value = '  +8.44%  '
Set-TextValue "D23" '13.17'
$ws.Range("E23").Value = '  +2.51%  '
Set-TextValue "D24" '27.356.92'
$ws.Range("E24").Value = '  +9.28%  '
Set-TextValue "D25" '2.464'
$ws.Range("E25").Value = '  +0.80%  '
Set-TextValue "D26" '1.544'
$ws.Range("E26").Value = '  +28.95%  '
Set-TextValue "D27" '2.559'
$ws.Range("E27").Value = '  +6.58%  '
Set-TextValue "D28" '20.36'
$ws.Range("E28").Value = '  +8.12%  '
Set-TextValue "D29" '153.88'
$ws.Range("E29").Value = '  +3.09%  '
Set-TextValue "D30" '1.987.13'
$ws.Range("E30").Value = '  +6.67%  '
Set-TextValue "D31" '135.65'
$ws.Range("E31").Value = '  +7.46%  '
Set-TextValue "D32" '6.402'
$ws.Range("E32").Value = '  +9.76%  '
Set-TextValue "D33" '4.182'
$ws.Range("E33").Value = '  +2.39%  '
Set-TextValue "D34" '0.08771'
$ws.Range("E34").Value = '  +3.50%  '
Set-TextValue "D35" '13.58'
$ws.Range("E35").Value = '  +9.44%  '
Set-TextValue "D36" '1.728'
$ws.Range("E36").Value = '  +3.48%  '
Set-TextValue "D37" '5.651'
$ws.Range("E37").Value = '  +8.97%  '
Set-TextValue "D38" '0.02427'
$ws.Range("E38").Value = '  +8.29%  '
$ws.Range("E39").Value = '  +8.57%  '
Set-TextValue "D40" '0.06534'
$ws.Range("E40").Value = '  +6.88%  '
Set-TextValue "D41" '0.6809'
$ws.Range("E41").Value = '  +13.46%  '
Set-TextValue "D42" '8.857'
$ws.Range("E42").Value = '  +6.31%  '
Set-TextValue "D43" '1.245'
$ws.Range("E43").Value = '  +1.00%  '
Set-TextValue "D44" '14.84'
$ws.Range("E44").Value = '  +8.66%  '
Set-TextValue "D45" '0.6450'
$ws.Range("E45").Value = '  +12.07%  '
Set-TextValue "D46" '0.9993'
$ws.Range("E46").Value = '  +0.09%  '
Set-TextValue "D47" '3.978'
$ws.Range("E47").Value = '  +3.66%  '
Set-TextValue "D48" '2.161'
$ws.Range("E48").Value = '  +9.74%  '
Set-TextValue "D49" '131.69'
$ws.Range("E49").Value = '  +5.76%  '
$ws.Range("E50").Value = '  +4.90%  '
Set-TextValue "D51" '80.59'
$ws.Range("E51").Value = '  +7.24%  '
